# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# with freshly scraped values, matching the source commit's data refresh.
# Values are written as literal text (Price/Volume columns store formatted,
# locale-style strings such as "29.899.65" or "  -0.31%  ", not numbers),
# so we force Text number formatting before assignment and then restore the
# default "Normal" style so no extra formatting is left behind.

function Set-CellText {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" '29.899.65'
Set-CellText $ws "E2" '  -0.31%  '
Set-CellText $ws "D3" '1.897.36'
Set-CellText $ws "E3" '  -0.02%  '
Set-CellText $ws "D4" '1.000'
Set-CellText $ws "E4" '  -0.01%  '
Set-CellText $ws "D5" '0.7933'
Set-CellText $ws "E5" '  -4.37%  '
Set-CellText $ws "D6" '244.00'
Set-CellText $ws "E6" '  +0.92%  '
Set-CellText $ws "E7" '  +0.02%  '
Set-CellText $ws "D8" '0.3164'
Set-CellText $ws "E8" '  -3.27%  '
Set-CellText $ws "D9" '25.44'
Set-CellText $ws "E9" '  -3.97%  '
Set-CellText $ws "D10" '0.07221'
Set-CellText $ws "E10" '  +2.76%  '
Set-CellText $ws "D11" '0.08106'
Set-CellText $ws "E11" '  +0.25%  '
Set-CellText $ws "E12" '  +0.80%  '
Set-CellText $ws "D13" '5.573'
Set-CellText $ws "E13" '  +6.29%  '
Set-CellText $ws "D14" '1.877.06'
Set-CellText $ws "E14" '  -1.21%  '
Set-CellText $ws "D15" '92.65'
Set-CellText $ws "E15" '  +0.57%  '
Set-CellText $ws "D16" '6.176'
Set-CellText $ws "E16" '  +5.78%  '
Set-CellText $ws "D17" '29.897.62'
Set-CellText $ws "E17" '  -0.32%  '
Set-CellText $ws "D18" '13.96'
Set-CellText $ws "E18" '  -0.86%  '
Set-CellText $ws "E19" '  +0.52%  '
Set-CellText $ws "D20" '0.000007803'
Set-CellText $ws "D21" '8.225'
Set-CellText $ws "E21" '  +18.45%  '
Set-CellText $ws "D22" '2.153.51'
Set-CellText $ws "E22" '  +0.22%  '
Set-CellText $ws "D23" '1.001'
Set-CellText $ws "E23" '  +0.01%  '
Set-CellText $ws "E24" '  +0.03%  '
Set-CellText $ws "D25" '0.1677'
Set-CellText $ws "E25" '  -3.63%  '
Set-CellText $ws "D26" '9.471'
Set-CellText $ws "E26" '  +2.38%  '
Set-CellText $ws "D27" '164.24'
Set-CellText $ws "E27" '  -0.68%  '
Set-CellText $ws "E28" '  -0.84%  '
Set-CellText $ws "D29" '2.065'
Set-CellText $ws "E29" '  -1.20%  '
Set-CellText $ws "D30" '1.399'
Set-CellText $ws "E30" '  +2.80%  '
Set-CellText $ws "E31" '  +2.26%  '
Set-CellText $ws "D32" '4.491'
Set-CellText $ws "E32" '  +5.11%  '
Set-CellText $ws "D33" '0.05574'
Set-CellText $ws "E33" '  -5.03%  '
Set-CellText $ws "D34" '4.098'
Set-CellText $ws "E34" '  +0.88%  '
Set-CellText $ws "D35" '1.283'
Set-CellText $ws "E35" '  +1.50%  '
Set-CellText $ws "D36" '0.7418'
Set-CellText $ws "E36" '  +1.58%  '
Set-CellText $ws "D37" '1.000'
Set-CellText $ws "E37" '  +0.02%  '
Set-CellText $ws "D38" '2.629'
Set-CellText $ws "E38" '  -3.29%  '
Set-CellText $ws "D39" '0.01933'
Set-CellText $ws "E39" '  +1.03%  '
Set-CellText $ws "D41" '1.152.38'
Set-CellText $ws "E41" '  +16.45%  '
Set-CellText $ws "D42" '74.53'
Set-CellText $ws "E42" '  +3.01%  '
Set-CellText $ws "E43" '  -0.17%  '
Set-CellText $ws "D44" '5.918'
Set-CellText $ws "E44" '  +1.19%  '
Set-CellText $ws "E45" '  +0.15%  '
Set-CellText $ws "D46" '104.91'
Set-CellText $ws "E46" '  +2.85%  '
Set-CellText $ws "D47" '1.001'
Set-CellText $ws "E47" '  +0.00%  '
Set-CellText $ws "D48" '1.884'
Set-CellText $ws "E48" '  -0.62%  '
Set-CellText $ws "D49" '10.04'
Set-CellText $ws "E49" '  +2.04%  '
Set-CellText $ws "D50" '3.042'
Set-CellText $ws "E50" '  +11.57%  '
Set-CellText $ws "D51" '7.466'
Set-CellText $ws "E51" '  -0.97%  '
